$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of H1 (header style) to I1/J1 so they match the existing header formatting
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 6
